$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the recalculated "Aggregate" biome min/RLI/max values (rows 2-37, years 1990-2025)
$aggUpdates = @(
  @(2, 0.95086826292113602, 0.95284026527353505, 0.95419883448951004),
  @(3, 0.95078169618543795, 0.95270824779020402, 0.95408251806023603),
  @(4, 0.95066048488452504, 0.95263175840875802, 0.954024317953307),
  @(5, 0.950622874783009, 0.95243869372326195, 0.95366522547566701),
  @(6, 0.95054905304437798, 0.95220067718901602, 0.95352398515217596),
  @(7, 0.95029471218160699, 0.95198916658676502, 0.95330657430681898),
  @(8, 0.950210149758279, 0.95181672633310599, 0.95304095243195297),
  @(9, 0.94999568076321095, 0.95162258618946904, 0.952814177582477),
  @(10, 0.94986237088532299, 0.95143277145728899, 0.95263961408878595),
  @(11, 0.949685697164381, 0.95123736387907898, 0.95238687319559301),
  @(12, 0.949580866691355, 0.95102015405503904, 0.95221806210098503),
  @(13, 0.94937320448368101, 0.95083251869353602, 0.95197279112077404),
  @(14, 0.94913088917246702, 0.95065137235338304, 0.95174488843958605),
  @(15, 0.94902718443059997, 0.95043203529248499, 0.95155083735380197),
  @(16, 0.94882497947837896, 0.95023970709652805, 0.95135800847883301),
  @(17, 0.94862321821656304, 0.950075376766225, 0.95118424569791504),
  @(18, 0.948291242845968, 0.94984896391454399, 0.95094834303439202),
  @(19, 0.94808472693406898, 0.94964111530342399, 0.95078097115160798),
  @(20, 0.94799804273230404, 0.94945235399199501, 0.95063740120566398),
  @(21, 0.94769937303601604, 0.94925521652145795, 0.95048068282305398),
  @(22, 0.94744296282079998, 0.94909712257914602, 0.95024468461746303),
  @(23, 0.94708882934767102, 0.94889332634224899, 0.950023517600326),
  @(24, 0.94696084647711198, 0.94873443344833697, 0.94992330973765604),
  @(25, 0.94669601223530697, 0.94853505597711496, 0.94974451701940699),
  @(26, 0.94633502691721305, 0.94833724058200797, 0.94958725121325105),
  @(27, 0.94616173653812197, 0.94808635069109803, 0.94948104557132995),
  @(28, 0.945946364430019, 0.94795816471676297, 0.94927226399420594),
  @(29, 0.94511017256630103, 0.94733173100933998, 0.948998460963145),
  @(30, 0.94308244274952402, 0.94639358079731295, 0.94853870725812595),
  @(31, 0.94084316144576896, 0.94498311811876401, 0.94812928959319298),
  @(32, 0.93890127720662298, 0.94264934695348201, 0.94705621261672401),
  @(33, 0.936143339602861, 0.94037708487034499, 0.94475901451715205),
  @(34, 0.93360021473550603, 0.93853894095923995, 0.94219561450865497),
  @(35, 0.93137986758134705, 0.93590642814234704, 0.94042647258373702),
  @(36, 0.93047599649977297, 0.93488832185706805, 0.93829708022972602),
  @(37, 0.92976825866919099, 0.93349892699743797, 0.93672353243611395)
)

foreach ($row in $aggUpdates) {
    $r = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]
}

# Add new "Forests" biome summary rows (65-67), matching the style of the
# preceding per-biome block (copy format from the last existing biome block)
$ws.Range("A62:F64").Copy()
[void]$ws.Range("A65:F67").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$forestRows = @(
  @(65, "Forests", 1990, 0.94727272727272704, 0.95272727272727298, 0.95272727272727298, 1990),
  @(66, "Forests", 2018, 0.941818181818182, 0.94727272727272704, 0.94727272727272704, 2018),
  @(67, "Forests", 2025, 0.94545454545454599, 0.94909090909090899, 0.94909090909090998, 2025)
)

foreach ($row in $forestRows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
}

# Move the selection/view to reflect where the new data was added
[void]$ws.Range("A67").Select()
